# A new daily price record (for "Provincia de Santiago", date 2022-06-13 /
# serial 44726) was inserted into the weekly "Apio" series at row 287,
# pushing all subsequent records (previously rows 287-404) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 287; this shifts the existing rows 287:404
# down to 288:405 and extends the used range to A1:R405.
$ws.Rows("287:287").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A287").Value = 3
$ws.Range("B287").Value = "Femacal de La Calera"
$ws.Range("C287").Value = "Coquimbo"
$ws.Range("D287").Value = 44726
$ws.Range("E287").Value = 5
$ws.Range("F287").Value = 100112017
$ws.Range("G287").Value = "Apio"
$ws.Range("H287").Value = "Americana (o)"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 230
$ws.Range("K287").Value = 8000
$ws.Range("L287").Value = 8500
$ws.Range("M287").Value = 8261
$ws.Range("N287").Value = "`$/docena de matas"
$ws.Range("O287").Value = "Provincia de Santiago"
$ws.Range("P287").Value = 1377
$ws.Range("Q287").Value = 6
$ws.Range("R287").Value = "Hortaliza"
